# Add a new LeetCode entry (row 38) to Sheet1, matching the upstream commit
# "Add files via upload" which appends the "Excel Sheet Column Number" problem.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row (row 38) -------------------------------------------------
# NOTE: the hyperlink is added to B38 *before* the cell text is written.
# Hyperlinks.Add() always stamps the cell with the supplied TextToDisplay
# (needed so the <hyperlink display="..."/> attribute is written out), so we
# let it run first and then overwrite the cell with the real problem name.
$ws.Hyperlinks.Add($ws.Range("B38"), "https://leetcode.com/problems/excel-sheet-column-number/", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/excel-sheet-column-number/")

$ws.Cells.Item(38, 2).Value2 = "Excel Sheet Column Number"
$ws.Cells.Item(38, 3).Value2 = 1
$ws.Cells.Item(38, 4).Value2 = 1
$ws.Cells.Item(38, 5).Value2 = 45
$ws.Cells.Item(38, 6).Value2 = 0.74
$ws.Cells.Item(38, 7).Value2 = 16.2
$ws.Cells.Item(38, 8).Value2 = 0.21
$ws.Cells.Item(38, 9).Value2 = "https://leetcode.com/problems/excel-sheet-column-number/submissions/1077070542/"

# Re-apply the same "Hyperlink" look used by the other problem-name cells.
# (Hyperlinks.Add() re-styles the cell with a brand new style record; copying
# the formatting of the previous hyperlink cell keeps it identical/shared.)
$ws.Range("B37").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection tweak (matches where Excel lands after the paste) -----------
$ws.Range("I38").Select() | Out-Null

# Recalculate so that Sheet2's summary formulas pick up the new row.
$excel.CalculateFull()
